$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume refresh: update the changed Price (D) and
# Volume(1h) (E) cells, plus the re-ranked Bittensor/Maker and
# ThetaToken/USDe rows (B/C/D/E). Price & volume cells are stored as
# text (e.g. "69.559.16", "0.999") in the source data, so force text
# formatting on exactly the cells being written before assigning the
# value - this avoids Excel auto-converting numeric-looking strings
# into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.559.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.734.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +8.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.05"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.731.21"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.85%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.63"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000253"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.345.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.724.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.592.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.91"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.41"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000125"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +19.71%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.19"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.337"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.40%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "422.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.087.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.76"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.51"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.01"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.17%  "
